# Atualizacao para APS Bacen - replace the sample row with the full
# "relacao de itens" (tintas / materiais de pintura) dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: Item Numero, Descricao, Quantidade, Valor Unitario (text), Valor Total
$rows = @(
    @{ Row = 2;  A = 90;  B = "TINTAANTIFERRUGE M"; C = 60; D = "152.0";  E = 9120 },
    @{ Row = 3;  A = 91;  B = "TINTA ESMALTE";       C = 30; D = "127.5";  E = 3825 },
    @{ Row = 4;  A = 92;  B = "TINTA ACRÍLICA";      C = 30; D = "308.0";  E = 9240 },
    @{ Row = 5;  A = 93;  B = "TINTA ACRÍLICA";      C = 20; D = "504.0";  E = 10080 },
    @{ Row = 6;  A = 94;  B = "TINTA CONCENTRADA";   C = 20; D = "319.0";  E = 6380 },
    @{ Row = 7;  A = 95;  B = "TINTA ACRÍLICA";      C = 15; D = "594.0";  E = 8910 },
    @{ Row = 8;  A = 96;  B = "DILUENTE";            C = 30; D = "118.0";  E = 3540 },
    @{ Row = 9;  A = 97;  B = "DILUENTE";            C = 20; D = "112.0";  E = 2240 },
    @{ Row = 10; A = 98;  B = "COLA";                C = 5;  D = "451.04"; E = 2255.2 },
    @{ Row = 11; A = 99;  B = "MASSA CORRIDA";       C = 30; D = "87.0";   E = 2610 },
    @{ Row = 12; A = 100; B = "MASSA CORRIDA";       C = 15; D = "172.0";  E = 2580 }
)

# The "Valor Unitario" column holds numeric-looking text (e.g. "152.0"); Excel
# auto-converts plain numeric strings to numbers on assignment, so format the
# whole column as Text first, then restore the default style afterwards so no
# stray style attribute is left on the cells (matches how the data was typed
# in as text in the source sheet).
$dColumn = $ws.Range("D2:D12")
$dColumn.NumberFormat = "@"

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 3).Value = $r.C
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}

$dColumn.Style = "Normal"
